# Plantilla de Casos de Uso - actualizacion de CU-20 (CRU profesor) y
# CU-21 (CRU cliente): se planifican ambos casos de uso y se carga su
# esfuerzo estimado (hrs).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# CU- 20 "CRU profesor" (fila 24): Estado vacio -> planificado, Esfuerzo 0 -> 2
$ws.Range("E24").Value = "planificado"
$ws.Range("F24").Value = 2

# CU- 21 "CRU cliente" (fila 25): Estado vacio -> planificado, Esfuerzo 0 -> 1
$ws.Range("E25").Value = "planificado"
$ws.Range("F25").Value = 1

# Deja la vista de la hoja tal como quedo guardada por el autor: desplazada
# para mostrar desde la fila 16 y con la celda E26 seleccionada.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("E26").Select()
